$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 910.7143
$ws.Range("I6").Value = 340
$ws.Range("J6").Value = 1227.7778
$ws.Range("K6").Value = 1020
$ws.Range("L6").Value = 3683.3334
$ws.Range("M6").Value = -908
$ws.Range("N6").Value = -3907.3334
$ws.Range("H8").Value = 123.5
$ws.Range("I8").Value = 123.5
$ws.Range("K8").Value = 370.5
$ws.Range("M8").Value = -231.5
$ws.Range("H12").Value = 84
$ws.Range("I12").Value = 84
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 84
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 86
$ws.Range("N12").ClearContents()
$ws.Range("H52").Value = 650
$ws.Range("I52").Value = 650
$ws.Range("K52").Value = 1950
$ws.Range("M52").Value = -1790
$ws.Range("H112").Value = 8523950
$ws.Range("J112").Value = 9405665
$ws.Range("L112").Value = 28216995
$ws.Range("N112").Value = -28219211
$ws.Range("H129").Value = 1064.9487
$ws.Range("J129").Value = 1182.1177
$ws.Range("L129").Value = 3546.3531
$ws.Range("N129").Value = -13546.3531
$ws.Range("H137").Value = 29413138
$ws.Range("I137").Value = 38462588
$ws.Range("J137").Value = 2425
$ws.Range("K137").Value = 115387764
$ws.Range("L137").Value = 7275
$ws.Range("M137").Value = -115385214
$ws.Range("N137").Value = -12375
$ws.Range("H138").Value = 9016262
$ws.Range("I138").Value = 2944342.5
$ws.Range("K138").Value = 8833027.5
$ws.Range("M138").Value = -8827887.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 34198.566
$ws.Range("I2").Value = 34198.566
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 34198.566
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -34085.566
$ws.Range("N2").ClearContents()
$ws.Range("H32").Value = 21095.973
$ws.Range("I32").Value = 4481.8867
$ws.Range("J32").Value = 67440.52
$ws.Range("K32").Value = 4481.8867
$ws.Range("L32").Value = 67440.52
$ws.Range("M32").Value = -4194.8867
$ws.Range("N32").Value = -68014.52
$ws.Range("H61").Value = 3247.6428
$ws.Range("I61").Value = 2569.9092
$ws.Range("J61").Value = 5732.6665
$ws.Range("K61").Value = 2569.9092
$ws.Range("L61").Value = 5732.6665
$ws.Range("M61").Value = -2357.9092
$ws.Range("N61").Value = -6156.6665
$ws.Range("H74").Value = 5657.9644
$ws.Range("I74").Value = 984.05
$ws.Range("J74").Value = 17342.75
$ws.Range("K74").Value = 984.05
$ws.Range("L74").Value = 17342.75
$ws.Range("M74").Value = -110.05
$ws.Range("N74").Value = -19090.75
$ws.Range("H77").Value = 5657.9644
$ws.Range("I77").Value = 984.05
$ws.Range("J77").Value = 17342.75
$ws.Range("K77").Value = 4920.25
$ws.Range("L77").Value = 86713.75
$ws.Range("M77").Value = -552.25
$ws.Range("N77").Value = -95449.75
$ws.Range("H116").Value = 34198.566
$ws.Range("I116").Value = 34198.566
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 34198.566
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -31904.566
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 3410.4827
$ws.Range("I132").Value = 2923
$ws.Range("K132").Value = 8769
$ws.Range("M132").Value = -6239
$ws.Range("H136").Value = 3247.6428
$ws.Range("I136").Value = 2569.9092
$ws.Range("J136").Value = 5732.6665
$ws.Range("K136").Value = 7709.7276
$ws.Range("L136").Value = 17197.9995
$ws.Range("M136").Value = -5159.7276
$ws.Range("N136").Value = -22297.9995
$ws.Range("H139").Value = 51882.375
$ws.Range("J139").Value = 51882.375
$ws.Range("L139").Value = 51882.375
$ws.Range("N139").Value = -62162.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 34198.566
$ws.Range("I3").Value = 34198.566
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 34198.566
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -34084.566
$ws.Range("N3").ClearContents()
$ws.Range("H134").Value = 3773.457
$ws.Range("I134").Value = 1858.3914
$ws.Range("J134").Value = 7444
$ws.Range("K134").Value = 5575.174199999999
$ws.Range("L134").Value = 22332
$ws.Range("M134").Value = -3040.174199999999
$ws.Range("N134").Value = -27402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6359.976
$ws.Range("I31").Value = 2146.1738
$ws.Range("J31").Value = 11460.895
$ws.Range("K31").Value = 2146.1738
$ws.Range("L31").Value = 11460.895
$ws.Range("M31").Value = -1851.1738
$ws.Range("N31").Value = -12050.895
$ws.Range("H34").Value = 6359.976
$ws.Range("I34").Value = 2146.1738
$ws.Range("J34").Value = 11460.895
$ws.Range("K34").Value = 2146.1738
$ws.Range("L34").Value = 11460.895
$ws.Range("M34").Value = -1944.1738
$ws.Range("N34").Value = -11864.895
$ws.Range("H36").Value = 2250
$ws.Range("I36").Value = 2250
$ws.Range("K36").Value = 2250
$ws.Range("M36").Value = -1862
$ws.Range("H40").Value = 2250
$ws.Range("I40").Value = 2250
$ws.Range("K40").Value = 2250
$ws.Range("M40").Value = -2090
$ws.Range("H99").Value = 2521.5557
$ws.Range("I99").Value = 1808
$ws.Range("J99").Value = 3642.8572
$ws.Range("K99").Value = 1808
$ws.Range("L99").Value = 3642.8572
$ws.Range("M99").Value = -310
$ws.Range("N99").Value = -6638.8572
$ws.Range("H105").Value = 1076.5
$ws.Range("I105").Value = 1074.4445
$ws.Range("J105").Value = 1080.2
$ws.Range("K105").Value = 1074.4445
$ws.Range("L105").Value = 1080.2
$ws.Range("M105").Value = 672.5554999999999
$ws.Range("N105").Value = -4574.2
$ws.Range("H126").Value = 2521.5557
$ws.Range("I126").Value = 1808
$ws.Range("J126").Value = 3642.8572
$ws.Range("K126").Value = 5424
$ws.Range("L126").Value = 10928.5716
$ws.Range("M126").Value = -2954
$ws.Range("N126").Value = -15868.5716
$ws.Range("H134").Value = 3271.72
$ws.Range("I134").Value = 1458
$ws.Range("J134").Value = 5580.091
$ws.Range("K134").Value = 4374
$ws.Range("L134").Value = 16740.273
$ws.Range("M134").Value = -1839
$ws.Range("N134").Value = -21810.273

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2467.923
$ws.Range("J34").Value = 2489
$ws.Range("L34").Value = 7467
$ws.Range("N34").Value = -7635
$ws.Range("H39").Value = 8406.666999999999
$ws.Range("J39").Value = 8711.538
$ws.Range("L39").Value = 26134.614
$ws.Range("N39").Value = -26722.614
$ws.Range("H107").Value = 384.5435
$ws.Range("J107").Value = 344.77274
$ws.Range("L107").Value = 1034.31822
$ws.Range("N107").Value = -4874.31822

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6120.5
$ws.Range("I22").Value = 855
$ws.Range("J22").Value = 9630.833000000001
$ws.Range("K22").Value = 855
$ws.Range("L22").Value = 9630.833000000001
$ws.Range("M22").Value = -560
$ws.Range("N22").Value = -10220.833
$ws.Range("H27").Value = 6120.5
$ws.Range("I27").Value = 855
$ws.Range("J27").Value = 9630.833000000001
$ws.Range("K27").Value = 855
$ws.Range("L27").Value = 9630.833000000001
$ws.Range("M27").Value = -748
$ws.Range("N27").Value = -9844.833000000001
$ws.Range("H46").Value = 1486.6666
$ws.Range("I46").Value = 1375
$ws.Range("J46").Value = 1527.2727
$ws.Range("K46").Value = 1375
$ws.Range("L46").Value = 1527.2727
$ws.Range("M46").Value = -1187
$ws.Range("N46").Value = -1903.2727
$ws.Range("H132").Value = 3423.6667
$ws.Range("I132").Value = 2049.087
$ws.Range("J132").Value = 5855.615
$ws.Range("K132").Value = 6147.261
$ws.Range("L132").Value = 17566.845
$ws.Range("M132").Value = -3617.261
$ws.Range("N132").Value = -22626.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13844.5
$ws.Range("J45").Value = 13844.5
$ws.Range("L45").Value = 13844.5
$ws.Range("N45").Value = -14826.5
$ws.Range("H81").Value = 2934.6943
$ws.Range("I81").Value = 1392.45
$ws.Range("J81").Value = 4862.5
$ws.Range("K81").Value = 2784.9
$ws.Range("L81").Value = 9725
$ws.Range("M81").Value = -1723.9
$ws.Range("N81").Value = -11847
$ws.Range("H84").Value = 2934.6943
$ws.Range("I84").Value = 1392.45
$ws.Range("J84").Value = 4862.5
$ws.Range("K84").Value = 13924.5
$ws.Range("L84").Value = 48625
$ws.Range("M84").Value = -8620.5
$ws.Range("N84").Value = -59233
$ws.Range("H94").Value = 27333.334
$ws.Range("J94").Value = 27333.334
$ws.Range("L94").Value = 27333.334
$ws.Range("N94").Value = -29135.334
$ws.Range("H126").Value = 48567
$ws.Range("I126").Value = 77454
$ws.Range("K126").Value = 232362
$ws.Range("M126").Value = -229892
